$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are numeric-looking text (e.g. "1.000", "23.476.75")
# that Excel would silently coerce to numbers on assignment. Forcing the
# cell to Text format before writing preserves the exact literal string,
# matching the inline-string cells in the source workbook.
$priceCells = @{
    "D2" = "23.476.75"
    "D3" = "1.633.24"
    "D4" = "1.000"
    "D5" = "1.001"
    "D6" = "305.60"
    "D7" = "0.3756"
    "D8" = "0.3672"
    "D9" = "51.76"
    "D10" = "0.08221"
    "D11" = "1.231"
    "D12" = "1.000"
    "D13" = "22.61"
    "D14" = "6.585"
    "D15" = "0.00001254"
    "D16" = "7.285"
    "D17" = "1.640.96"
    "D18" = "94.31"
    "D19" = "0.06980"
    "D20" = "17.82"
    "D21" = "6.472"
    "D22" = "1.000"
    "D23" = "12.79"
    "D24" = "23.480.70"
    "D25" = "3.177"
    "D26" = "2.462"
    "D27" = "21.45"
    "D28" = "150.43"
    "D29" = "5.332"
    "D30" = "134.60"
    "D31" = "1.818.80"
    "D32" = "2.277"
    "D33" = "6.841"
    "D34" = "1.029"
    "D35" = "11.00"
    "D36" = "0.02799"
    "D37" = "0.2539"
    "D38" = "6.101"
    "D39" = "0.07191"
    "D40" = "0.08772"
    "D41" = "0.7093"
    "D42" = "1.353"
    "D43" = "16.34"
    "D44" = "12.32"
    "D45" = "0.6588"
    "D46" = "2.342"
    "D47" = "0.9999"
    "D48" = "4.002"
    "D49" = "0.08029"
    "D50" = "1.211"
    "D51" = "125.85"
}
foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
}

# Volume(1h) column (E) values already contain spaces/%% so Excel keeps them as text.
$volumeCells = @{
    "E2" = "  -0.36%  "
    "E3" = "  -0.49%  "
    "E4" = "  +0.02%  "
    "E5" = "  +0.03%  "
    "E6" = "  -1.04%  "
    "E7" = "  -0.37%  "
    "E8" = "  +0.12%  "
    "E9" = "  -1.54%  "
    "E10" = "  +0.17%  "
    "E11" = "  -3.68%  "
    "E12" = "  +0.01%  "
    "E13" = "  -1.86%  "
    "E14" = "  -1.25%  "
    "E15" = "  -2.54%  "
    "E16" = "  -2.00%  "
    "E17" = "  +0.00%  "
    "E18" = "  -0.60%  "
    "E19" = "  +0.62%  "
    "E20" = "  -2.64%  "
    "E21" = "  -1.61%  "
    "E22" = "  +0.17%  "
    "E23" = "  -0.57%  "
    "E24" = "  -0.34%  "
    "E25" = "  +3.45%  "
    "E26" = "  +1.57%  "
    "E27" = "  +0.47%  "
    "E28" = "  -0.70%  "
    "E29" = "  -0.78%  "
    "E30" = "  -0.93%  "
    "E31" = "  -0.33%  "
    "E32" = "  -4.83%  "
    "E33" = "  +0.25%  "
    "E34" = "  +5.37%  "
    "E35" = "  +5.81%  "
    "E36" = "  -1.56%  "
    "E37" = "  -0.70%  "
    "E38" = "  -1.59%  "
    "E39" = "  -2.70%  "
    "E40" = "  -1.59%  "
    "E41" = "  -0.58%  "
    "E42" = "  -2.40%  "
    "E43" = "  +0.09%  "
    "E44" = "  -1.81%  "
    "E45" = "  +0.27%  "
    "E46" = "  -0.44%  "
    "E47" = "  +0.06%  "
    "E48" = "  -1.11%  "
    "E49" = "  +0.52%  "
    "E50" = "  -0.40%  "
    "E51" = "  -3.26%  "
}
foreach ($addr in $volumeCells.Keys) {
    $ws.Range($addr).Value = $volumeCells[$addr]
}

# Rows 38-40 also swapped which coin occupies which rank (Coin name + Link).
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
